$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.09912404012357631

$ws.Range("B3").Value = 0.001969907350911
$ws.Range("C3").Value = 0.0006792740834497284
$ws.Range("D3").Value = 3.989874991484508
$ws.Range("E3").Value = 0.08790630064225226
$ws.Range("F3").Value = 0.000638546645195661
$ws.Range("G3").Value = 0.003301268056626339
$ws.Range("H3").Value = 0.1010939474744873

$ws.Range("B4").Value = 0.004763016599688532
$ws.Range("C4").Value = 0.001087764048200155
$ws.Range("D4").Value = 4.155518924700853
$ws.Range("E4").Value = 0.005041401044172722
$ws.Range("F4").Value = 0.002631027440920148
$ws.Range("G4").Value = 0.006895005758456918
$ws.Range("H4").Value = 0.1038870567232648

$ws.Range("B5").Value = 0.007961018694175597
$ws.Range("C5").Value = 0.001696903254367539
$ws.Range("D5").Value = 6.247422250883083
$ws.Range("E5").Value = 0.1893882024200223
$ws.Range("F5").Value = 0.004635133330843471
$ws.Range("G5").Value = 0.01128690405750772
$ws.Range("H5").Value = 0.1070850588177519

$ws.Range("B6").Value = 0.007713261819120072
$ws.Range("C6").Value = 0.005582290229956171
$ws.Range("D6").Value = 3.413103656481359
$ws.Range("E6").Value = 0.1275911100851306
$ws.Range("F6").Value = -0.003227886564427699
$ws.Range("G6").Value = 0.01865441020266784
$ws.Range("H6").Value = 0.1068373019426964

$ws.Range("B7").Value = 0.01248176075846321
$ws.Range("C7").Value = 0.004050449129408803
$ws.Range("D7").Value = 4.277125696196287
$ws.Range("E7").Value = 0.02440828912289726
$ws.Range("F7").Value = 0.004543002325326651
$ws.Range("G7").Value = 0.02042051919159977
$ws.Range("H7").Value = 0.1116058008820395

$ws.Range("B8").Value = 0.007873373587688241
$ws.Range("C8").Value = 0.001960956588305206
$ws.Range("D8").Value = 5.3650575923961
$ws.Range("E8").Value = 0.007060969426136016
$ws.Range("F8").Value = 0.004029952485432458
$ws.Range("G8").Value = 0.01171679468994402
$ws.Range("H8").Value = 0.1069974137112646

$ws.Range("B9").Value = 0.06126548629181352
$ws.Range("C9").Value = 0.005974746490127017
$ws.Range("D9").Value = 6.254955807264902
$ws.Range("E9").Value = 0.0002807193490298925
$ws.Range("F9").Value = 0.04955513648773151
$ws.Range("G9").Value = 0.07297583609589553
$ws.Range("H9").Value = 0.1603895264153898

$ws.Range("B10").Value = -0.09912404012357631
$ws.Range("C10").Value = 0.0005522708947199429
$ws.Range("D10").Value = -215.9199853302586
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.1002064776500371
$ws.Range("G10").Value = -0.09804160259711552
$ws.Range("H10").Value = 0

$ws.Range("B11").Value = -0.04429857794053839
$ws.Range("C11").Value = 0.0005979279000378096
$ws.Range("D11").Value = -86.6089600510645
$ws.Range("E11").Value = [double]"2.558580134952391e-261"
$ws.Range("F11").Value = -0.04547050198677859
$ws.Range("G11").Value = -0.04312665389429819
$ws.Range("H11").Value = 0.05482546218303792

$ws.Range("B12").Value = -0.0382895251006394
$ws.Range("C12").Value = 0.0005965268347801323
$ws.Range("D12").Value = -73.73003104183842
$ws.Range("E12").Value = [double]"2.045521105960129e-206"
$ws.Range("F12").Value = -0.03945870311004777
$ws.Range("G12").Value = -0.03712034709123102
$ws.Range("H12").Value = 0.06083451502293691

$ws.Range("B13").Value = -0.03667023456129203
$ws.Range("C13").Value = 0.0005950305101288869
$ws.Range("D13").Value = -69.26452122847117
$ws.Range("E13").Value = [double]"2.047807876872883e-172"
$ws.Range("F13").Value = -0.03783647988304711
$ws.Range("G13").Value = -0.03550398923953695
$ws.Range("H13").Value = 0.06245380556228428

$ws.Range("B14").Value = -0.03046362672923095
$ws.Range("C14").Value = 0.0005919712418413596
$ws.Range("D14").Value = -55.50374299386527
$ws.Range("E14").Value = [double]"4.992734507483147e-38"
$ws.Range("F14").Value = -0.03162387580609818
$ws.Range("G14").Value = -0.02930337765236372
$ws.Range("H14").Value = 0.06866041339434537

$ws.Range("B15").Value = -0.02764508677365244
$ws.Range("C15").Value = 0.0005762300664488778
$ws.Range("D15").Value = -52.69839258492829
$ws.Range("E15").Value = [double]"1.25033539213733e-27"
$ws.Range("F15").Value = -0.02877448362386029
$ws.Range("G15").Value = -0.02651568992344459
$ws.Range("H15").Value = 0.07147895334992388

$ws.Range("B16").Value = -0.02265940967752368
$ws.Range("C16").Value = 0.0005616018022999292
$ws.Range("D16").Value = -44.21827596319378
$ws.Range("E16").Value = [double]"2.074978406278525e-07"
$ws.Range("F16").Value = -0.02376013567352311
$ws.Range("G16").Value = -0.02155868368152426
$ws.Range("H16").Value = 0.07646463044605263

$ws.Range("B17").Value = -0.01835416086881763
$ws.Range("C17").Value = 0.0005650675835742472
$ws.Range("D17").Value = -38.47043005310683
$ws.Range("E17").Value = 0.04287914853338819
$ws.Range("F17").Value = -0.01946167965049072
$ws.Range("G17").Value = -0.01724664208714455
$ws.Range("H17").Value = 0.08076987925475868

$ws.Range("B18").Value = -0.01853491288746079
$ws.Range("C18").Value = 0.0005715171371309402
$ws.Range("D18").Value = -37.51008626564604
$ws.Range("E18").Value = [double]"2.846794667548497e-14"
$ws.Range("F18").Value = -0.01965507259681564
$ws.Range("G18").Value = -0.01741475317810594
$ws.Range("H18").Value = 0.08058912723611553

$ws.Range("B19").Value = -0.01617550511239547
$ws.Range("C19").Value = 0.0006019181169608974
$ws.Range("D19").Value = -32.42528235001339
$ws.Range("E19").Value = [double]"4.872812294569726e-10"
$ws.Range("F19").Value = -0.01735525066932838
$ws.Range("G19").Value = -0.01499575955546255
$ws.Range("H19").Value = 0.08294853501118085

$ws.Range("B20").Value = -0.01234987996562419
$ws.Range("C20").Value = 0.0006130558544811643
$ws.Range("D20").Value = -25.94470497262839
$ws.Range("E20").Value = 0.05275932284806881
$ws.Range("F20").Value = -0.01355145546544801
$ws.Range("G20").Value = -0.01114830446580037
$ws.Range("H20").Value = 0.08677416015795213

$ws.Range("B21").Value = -0.01174631605031217
$ws.Range("C21").Value = 0.0006213063200913889
$ws.Range("D21").Value = -21.25124396978418
$ws.Range("E21").Value = 0.01954117352813525
$ws.Range("F21").Value = -0.01296406212643655
$ws.Range("G21").Value = -0.0105285699741878
$ws.Range("H21").Value = 0.08737772407326413

$ws.Range("B22").Value = -0.009784054758511131
$ws.Range("C22").Value = 0.0005915846990609998
$ws.Range("D22").Value = -17.26145967619875
$ws.Range("E22").Value = 0.03498361640396983
$ws.Range("F22").Value = -0.0109435465971688
$ws.Range("G22").Value = -0.008624562919853465
$ws.Range("H22").Value = 0.08933998536506518

$ws.Range("B23").Value = -0.006674684167149318
$ws.Range("C23").Value = 0.0005895384108173295
$ws.Range("D23").Value = -12.71002665266026
$ws.Range("E23").Value = 0.01613324950104657
$ws.Range("F23").Value = -0.007830165407939172
$ws.Range("G23").Value = -0.005519202926359465
$ws.Range("H23").Value = 0.09244935595642699

$ws.Range("B24").Value = -0.003862601735944289
$ws.Range("C24").Value = 0.000578434674261406
$ws.Range("D24").Value = -7.55886471062492
$ws.Range("E24").Value = 0.0513026462280677
$ws.Range("F24").Value = -0.004996319809467696
$ws.Range("G24").Value = -0.00272888366242088
$ws.Range("H24").Value = 0.09526143838763203

$ws.Range("B25").Value = -0.002986180309636185
$ws.Range("C25").Value = 0.0005482076431902245
$ws.Range("D25").Value = -6.909494540869561
$ws.Range("E25").Value = 0.05400323717021923
$ws.Range("F25").Value = -0.004060654014139587
$ws.Range("G25").Value = -0.001911706605132783
$ws.Range("H25").Value = 0.09613785981394013

$ws.Range("B26").Value = 0.02834203878765998
$ws.Range("C26").Value = 0.0008328571950808638
$ws.Range("D26").Value = 23.98614778046045
$ws.Range("E26").Value = [double]"5.981163700908449e-05"
$ws.Range("F26").Value = 0.02670965929688141
$ws.Range("G26").Value = 0.02997441827843855
$ws.Range("H26").Value = 0.1274660789112363
